$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "64.129.62"
$ws.Range("E2").Value = "  +2.00%  "
Set-TextValue "D3" "3.144.12"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  -0.02%  "
Set-TextValue "D5" "592.06"
$ws.Range("E5").Value = "  +0.80%  "
Set-TextValue "D6" "147.47"
$ws.Range("E6").Value = "  +2.00%  "
$ws.Range("E7").Value = "  +0.06%  "
Set-TextValue "D8" "3.138.26"
$ws.Range("E8").Value = "  +1.16%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("E10").Value = "  +12.66%  "
Set-TextValue "D11" "5.75"
$ws.Range("E11").Value = "  -0.28%  "
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("E13").Value = "  +4.42%  "
Set-TextValue "D14" "37.54"
$ws.Range("E15").Value = "  -0.88%  "
Set-TextValue "D16" "3.664.49"
$ws.Range("E16").Value = "  +1.17%  "
Set-TextValue "D17" "63.990.63"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("E18").Value = "  -1.72%  "
Set-TextValue "D19" "3.140.66"
$ws.Range("E19").Value = "  +1.32%  "
Set-TextValue "D20" "468.44"
$ws.Range("E20").Value = "  +3.27%  "
Set-TextValue "D21" "14.39"
$ws.Range("E21").Value = "  +1.65%  "
Set-TextValue "D22" "0.738"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("E23").Value = "  +0.65%  "
Set-TextValue "D24" "13.33"
$ws.Range("E24").Value = "  -3.49%  "
Set-TextValue "D25" "82.61"
$ws.Range("E25").Value = "  +0.23%  "
$ws.Range("E26").Value = "  -0.20%  "
Set-TextValue "D27" "9.01"
$ws.Range("E27").Value = "  +8.46%  "
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("E29").Value = "  -1.52%  "
$ws.Range("E30").Value = "  +0.04%  "
Set-TextValue "D31" "6.88"
$ws.Range("E31").Value = "  -0.02%  "
Set-TextValue "D32" "27.24"
$ws.Range("E32").Value = "  +0.07%  "
$ws.Range("E33").Value = "  -4.74%  "
Set-TextValue "D34" "0.0₃0880"
$ws.Range("E34").Value = "  +9.59%  "
Set-TextValue "D35" "2.39"
$ws.Range("E35").Value = "  +7.64%  "
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("E37").Value = "  +13.65%  "
$ws.Range("E38").Value = "  +0.82%  "
Set-TextValue "D39" "461.08"
$ws.Range("E39").Value = "  +8.11%  "
Set-TextValue "D40" "51.01"
$ws.Range("E40").Value = "  +0.87%  "
Set-TextValue "D41" "8.75"
$ws.Range("E41").Value = "  -1.32%  "
Set-TextValue "D42" "0.0374"
$ws.Range("E42").Value = "  +0.81%  "
Set-TextValue "D43" "2.908.08"
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("E45").Value = "  +1.85%  "
$ws.Range("E46").Value = "  +0.79%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  +0.06%  "
Set-TextValue "D51" "24.87"
$ws.Range("E51").Value = "  +0.55%  "

# Row 47/48: Arweave and Monero swap places with updated data
$ws.Range("B47").Value = "Arweave"
$ws.Range("C47").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue "D47" "35.86"
$ws.Range("E47").Value = "  +2.16%  "

$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D48" "127.02"
$ws.Range("E48").Value = "  +1.85%  "
